# Update "想去人数" (interested-count) figures across sheets to match the
# newly generated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 223
$ws1.Range("F6").Value  = 808
$ws1.Range("F8").Value  = 10314
$ws1.Range("F12").Value = 2460
$ws1.Range("F15").Value = 1395
$ws1.Range("F21").Value = 397
$ws1.Range("F26").Value = 240
$ws1.Range("F34").Value = 3875
$ws1.Range("F35").Value = 3278
$ws1.Range("F36").Value = 35
$ws1.Range("F38").Value = 1049
$ws1.Range("F39").Value = 408

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 2079

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 808
$ws4.Range("F11").Value = 10315
$ws4.Range("F15").Value = 2460
$ws4.Range("F23").Value = 397
$ws4.Range("F26").Value = 240
$ws4.Range("F36").Value = 3278
$ws4.Range("F37").Value = 1049
